# Refresh the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with the
# latest symbol-list snapshot (GitHub Actions refresh run).
# Values are entered with a leading apostrophe so Excel stores them as the
# literal text strings used by the source feed (e.g. "4.50%"), matching the
# existing inline-string cells instead of coercing to Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'326.73"
$ws.Range("E2").Value = "'4.50%"
$ws.Range("D3").Value = "'39.73"
$ws.Range("E3").Value = "'5.89%"
$ws.Range("D4").Value = "'5.247"
$ws.Range("E4").Value = "'1.99%"
$ws.Range("D5").Value = "'0.08104"
$ws.Range("E5").Value = "'2.47%"
$ws.Range("E6").Value = "'2.36%"
$ws.Range("E7").Value = "'4.48%"
$ws.Range("D8").Value = "'1.920"
$ws.Range("E8").Value = "'0.27%"
$ws.Range("D10").Value = "'0.9342"
$ws.Range("E10").Value = "'0.56%"
$ws.Range("D11").Value = "'0.1319"
$ws.Range("E11").Value = "'19.08%"
$ws.Range("E12").Value = "'3.33%"
$ws.Range("D13").Value = "'0.09274"
$ws.Range("E13").Value = "'1.75%"
$ws.Range("D14").Value = "'0.03437"
$ws.Range("E14").Value = "'3.33%"
$ws.Range("D15").Value = "'0.09550"
$ws.Range("E15").Value = "'-0.70%"
$ws.Range("D16").Value = "'0.001394"
$ws.Range("E16").Value = "'0.64%"
$ws.Range("D17").Value = "'0.04443"
$ws.Range("E17").Value = "'1.60%"
$ws.Range("D18").Value = "'0.006134"
$ws.Range("E18").Value = "'7.50%"
$ws.Range("D19").Value = "'3.359"
$ws.Range("D20").Value = "'0.3535"
$ws.Range("E20").Value = "'3.76%"
$ws.Range("D21").Value = "'7.233"
$ws.Range("E21").Value = "'21.96%"
$ws.Range("D22").Value = "'0.1322"
$ws.Range("E22").Value = "'2.66%"
$ws.Range("E23").Value = "'-10.85%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-1.22%"
$ws.Range("D25").Value = "'0.004361"
$ws.Range("E25").Value = "'-5.98%"
$ws.Range("E26").Value = "'-5.17%"
$ws.Range("E27").Value = "'-0.13%"
$ws.Range("D39").Value = "'0.02489"
$ws.Range("E39").Value = "'10.67%"
$ws.Range("D40").Value = "'0.05239"
$ws.Range("E40").Value = "'2.94%"
$ws.Range("D41").Value = "'0.007687"
$ws.Range("E41").Value = "'3.11%"
$ws.Range("D42").Value = "'0.1432"
$ws.Range("E42").Value = "'5.78%"
$ws.Range("D43").Value = "'0.008603"
$ws.Range("E43").Value = "'-4.66%"
$ws.Range("E44").Value = "'-0.96%"
$ws.Range("D45").Value = "'0.008169"
$ws.Range("E45").Value = "'-5.51%"
$ws.Range("D46").Value = "'0.00006657"
$ws.Range("E46").Value = "'-0.60%"
$ws.Range("D48").Value = "'0.002852"
$ws.Range("E48").Value = "'-13.14%"
$ws.Range("E49").Value = "'148.02%"
